# Updated cryptos list — apply latest price/volume snapshot to Sheet1.
# All target cells are plain text cells (t="inlineStr" in the source OOXML),
# so every write forces text storage (NumberFormat "@" then ClearFormats)
# to stop Excel's autoconvert from turning numeric-looking strings (e.g.
# "605.49", "0.0000183") into real numbers, while leaving the cell's
# original (unstyled) formatting untouched afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Ref = 'D2'; Value = '70.933.56' },
    @{ Ref = 'E2'; Value = '  +4.75%  ' },
    @{ Ref = 'D3'; Value = '2.606.44' },
    @{ Ref = 'E3'; Value = '  +5.36%  ' },
    @{ Ref = 'E4'; Value = '  -0.02%  ' },
    @{ Ref = 'D5'; Value = '605.49' },
    @{ Ref = 'E5'; Value = '  +3.23%  ' },
    @{ Ref = 'D6'; Value = '181.45' },
    @{ Ref = 'E6'; Value = '  +3.49%  ' },
    @{ Ref = 'E7'; Value = '  -0.04%  ' },
    @{ Ref = 'E8'; Value = '  +1.95%  ' },
    @{ Ref = 'D9'; Value = '2.606.10' },
    @{ Ref = 'E10'; Value = '  +16.75%  ' },
    @{ Ref = 'E11'; Value = '  +0.53%  ' },
    @{ Ref = 'E12'; Value = '  +4.78%  ' },
    @{ Ref = 'E13'; Value = '  +1.94%  ' },
    @{ Ref = 'D14'; Value = '26.82' },
    @{ Ref = 'E14'; Value = '  +6.40%  ' },
    @{ Ref = 'B15'; Value = 'ShibaInu' },
    @{ Ref = 'C15'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' },
    @{ Ref = 'D15'; Value = '0.0000183' },
    @{ Ref = 'E15'; Value = '  +8.18%  ' },
    @{ Ref = 'B16'; Value = 'WrappedBTC' },
    @{ Ref = 'C16'; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc' },
    @{ Ref = 'D16'; Value = '70.961.49' },
    @{ Ref = 'E16'; Value = '  +4.80%  ' },
    @{ Ref = 'B17'; Value = 'WrappedliquidstakedEther2.0' },
    @{ Ref = 'C17'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth' },
    @{ Ref = 'D17'; Value = '2.982.62' },
    @{ Ref = 'E17'; Value = '  +1.97%  ' },
    @{ Ref = 'D18'; Value = '2.620.32' },
    @{ Ref = 'E18'; Value = '  +5.28%  ' },
    @{ Ref = 'B19'; Value = 'BitcoinCash' },
    @{ Ref = 'C19'; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch' },
    @{ Ref = 'D19'; Value = '378.92' },
    @{ Ref = 'E19'; Value = '  +8.95%  ' },
    @{ Ref = 'B20'; Value = 'Uniswap' },
    @{ Ref = 'C20'; Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni' },
    @{ Ref = 'D20'; Value = '7.92' },
    @{ Ref = 'E20'; Value = '  +7.21%  ' },
    @{ Ref = 'E21'; Value = '  +5.72%  ' },
    @{ Ref = 'E22'; Value = '  +3.66%  ' },
    @{ Ref = 'D23'; Value = '72.06' },
    @{ Ref = 'E23'; Value = '  +1.99%  ' },
    @{ Ref = 'E24'; Value = '  +5.44%  ' },
    @{ Ref = 'E25'; Value = '  +0.28%  ' },
    @{ Ref = 'E26'; Value = '  +12.48%  ' },
    @{ Ref = 'D27'; Value = '9.71' },
    @{ Ref = 'E27'; Value = '  +10.62%  ' },
    @{ Ref = 'D28'; Value = '2.738.64' },
    @{ Ref = 'E28'; Value = '  +5.73%  ' },
    @{ Ref = 'E29'; Value = '  +0.00%  ' },
    @{ Ref = 'D30'; Value = '0.0₃0949' },
    @{ Ref = 'E30'; Value = '  +7.07%  ' },
    @{ Ref = 'D31'; Value = '530.26' },
    @{ Ref = 'E31'; Value = '  +7.68%  ' },
    @{ Ref = 'D32'; Value = '8.08' },
    @{ Ref = 'E32'; Value = '  +5.11%  ' },
    @{ Ref = 'E33'; Value = '  +6.94%  ' },
    @{ Ref = 'E34'; Value = '  +4.80%  ' },
    @{ Ref = 'E35'; Value = '  -0.01%  ' },
    @{ Ref = 'D36'; Value = '163.98' },
    @{ Ref = 'E36'; Value = '  +0.59%  ' },
    @{ Ref = 'E37'; Value = '  +1.42%  ' },
    @{ Ref = 'E38'; Value = '  +5.19%  ' },
    @{ Ref = 'D39'; Value = '18.92' },
    @{ Ref = 'E39'; Value = '  +1.57%  ' },
    @{ Ref = 'E40'; Value = '  +7.15%  ' },
    @{ Ref = 'E41'; Value = '  +6.98%  ' },
    @{ Ref = 'E42'; Value = '  +6.38%  ' },
    @{ Ref = 'E43'; Value = '  +0.11%  ' },
    @{ Ref = 'E44'; Value = '  +7.05%  ' },
    @{ Ref = 'E45'; Value = '  +1.88%  ' },
    @{ Ref = 'D46'; Value = '40.50' },
    @{ Ref = 'E46'; Value = '  +4.61%  ' },
    @{ Ref = 'D47'; Value = '154.03' },
    @{ Ref = 'E47'; Value = '  +4.89%  ' },
    @{ Ref = 'E48'; Value = '  +4.36%  ' },
    @{ Ref = 'E49'; Value = '  +6.68%  ' },
    @{ Ref = 'D50'; Value = '0.533' },
    @{ Ref = 'E50'; Value = '  +4.82%  ' },
    @{ Ref = 'E51'; Value = '  +7.49%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}

Write-Host "Applied $($updates.Count) cell updates"
